# Apply the changes described by the diff.
$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 0.3.0 -> 0.4.0-snapshot-1
$meta.Range("B3").Value = "0.4.0-snapshot-1"

# Status: active -> draft
$meta.Range("B6").Value = "draft"

# Date: 2024-03-13T09:33:00+00:00 -> 2024-05-23T12:16:26+00:00
$meta.Range("B8").Value = "2024-05-23T12:16:26+00:00"

# Contact: "No display for ContactDetail" -> "ANS (https://esante.gouv.fr)"
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Swap the two "Mapping" header columns (AK <-> AL) and their widths/content
# so that the "Spécification métier" mapping now comes before "RIM Mapping".
$headerAK = $elements.Range("AK1").Value2
$headerAL = $elements.Range("AL1").Value2
$elements.Range("AK1").Value = $headerAL
$elements.Range("AL1").Value = $headerAK

$elements.Columns.Item(37).ColumnWidth = 78.0
$elements.Columns.Item(38).ColumnWidth = 24.166666666666668

# Swap the per-row mapping values to stay aligned with the swapped headers.
for ($r = 2; $r -le 6; $r++) {
    $ak = $elements.Cells.Item($r, 37).Value2
    $al = $elements.Cells.Item($r, 38).Value2
    $elements.Cells.Item($r, 37).Value = $al
    $elements.Cells.Item($r, 38).Value = $ak
}
